$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 59 - SEC v. Middleton, et al. (Veritaseum)
$ws.Range("D59").Value = "Ongoing"
$ws.Range("E59").Value = "Unregistered Offering and Fraud"
$ws.Range("F59").Value = "Civil"
$ws.Range("G59").Value = "VERI"
$ws.Range("H59").Value = "Veritaseum, Inc. and Veritaseum, LLC"
$ws.Range("I59").Value = "Ethereum"
$ws.Range("J59").Value = 14800000
$ws.Range("K59").Value = 1
$ws.Range("L59").Value = 1
$ws.Range("M59").Value = "New York"

# Row 60 - SEC v. Longfin Corp., et al.
$ws.Range("D60").Value = "Settlement"
$ws.Range("E60").Value = "Unregistered Offering"
$ws.Range("F60").Value = "Civil"
$ws.Range("G60").Value = "N/A"
$ws.Range("H60").Value = "Longfin Corp."
$ws.Range("I60").Value = "Ethereum"
$ws.Range("J60").Value = 27000000
$ws.Range("K60").Value = 1
$ws.Range("L60").Value = 1
$ws.Range("M60").Value = "Washington, D.C."

# Row 61 - SEC v. Longfin Corp., et al. (fraud claims)
$ws.Range("D61").Value = "Settlement"
$ws.Range("E61").Value = "Fraud"
$ws.Range("F61").Value = "Civil"
$ws.Range("G61").Value = "N/A"
$ws.Range("H61").Value = "Longfin Corp."
$ws.Range("I61").Value = "Ethereum"
$ws.Range("J61").Value = 27000000
$ws.Range("K61").Value = 1
$ws.Range("L61").Value = 1
$ws.Range("M61").Value = "Washington, D.C."

# Update view state to match the saved workbook: new active selection
$ws.Range("J62").Select() | Out-Null
